function Set-TextValue($ws, $addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "41.941.16"
Set-TextValue $ws "E2" "  -1.22%  "
Set-TextValue $ws "D3" "2.210.76"
Set-TextValue $ws "E3" "  -1.56%  "
Set-TextValue $ws "E4" "  +0.13%  "
Set-TextValue $ws "D5" "240.84"
Set-TextValue $ws "E5" "  -2.30%  "
Set-TextValue $ws "E6" "  -1.14%  "
Set-TextValue $ws "D7" "72.09"
Set-TextValue $ws "E7" "  -5.36%  "
Set-TextValue $ws "E8" "  +0.12%  "
Set-TextValue $ws "D9" "0.599"
Set-TextValue $ws "E9" "  -3.95%  "
Set-TextValue $ws "D10" "41.69"
Set-TextValue $ws "E10" "  -5.67%  "
Set-TextValue $ws "E11" "  +0.31%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D12" "0.104"
Set-TextValue $ws "E12" "  +0.76%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D13" "6.95"
Set-TextValue $ws "E13" "  -4.89%  "
Set-TextValue $ws "D14" "2.545.24"
Set-TextValue $ws "E14" "  -1.49%  "
Set-TextValue $ws "D15" "14.15"
Set-TextValue $ws "E15" "  -3.24%  "
Set-TextValue $ws "D16" "0.831"
Set-TextValue $ws "E16" "  -3.21%  "
Set-TextValue $ws "D17" "2.200.96"
Set-TextValue $ws "E17" "  -2.14%  "
Set-TextValue $ws "D18" "41.781.67"
Set-TextValue $ws "E18" "  -1.03%  "
Set-TextValue $ws "E19" "  +2.97%  "
Set-TextValue $ws "E20" "  +0.17%  "
Set-TextValue $ws "D21" "6.12"
Set-TextValue $ws "E21" "  -1.09%  "
Set-TextValue $ws "D22" "10.93"
Set-TextValue $ws "E22" "  +19.64%  "
Set-TextValue $ws "D23" "229.01"
Set-TextValue $ws "E23" "  -1.16%  "
Set-TextValue $ws "E24" "  -8.71%  "
Set-TextValue $ws "D25" "1.00"
Set-TextValue $ws "E25" "  -0.04%  "
Set-TextValue $ws "E26" "  -0.27%  "
Set-TextValue $ws "E27" "  +0.13%  "
Set-TextValue $ws "E28" "  -2.64%  "
Set-TextValue $ws "E29" "  -0.81%  "
Set-TextValue $ws "D30" "167.32"
Set-TextValue $ws "E30" "  -0.56%  "
Set-TextValue $ws "D31" "20.37"
Set-TextValue $ws "E31" "  -1.63%  "
Set-TextValue $ws "D32" "5.60"
Set-TextValue $ws "E32" "  +6.08%  "
Set-TextValue $ws "D33" "0.0792"
Set-TextValue $ws "E33" "  -4.65%  "
Set-TextValue $ws "D34" "29.90"
Set-TextValue $ws "E34" "  -2.70%  "
Set-TextValue $ws "E35" "  -1.26%  "
Set-TextValue $ws "E36" "  -12.32%  "
Set-TextValue $ws "E37" "  -7.86%  "
Set-TextValue $ws "D38" "0.0299"
Set-TextValue $ws "E38" "  -6.29%  "
Set-TextValue $ws "D39" "13.64"
Set-TextValue $ws "E39" "  -0.79%  "
Set-TextValue $ws "E40" "  -4.01%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws "D41" "63.90"
Set-TextValue $ws "E41" "  +0.37%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws "D42" "5.59"
Set-TextValue $ws "E42" "  -3.82%  "
Set-TextValue $ws "D43" "0.195"
Set-TextValue $ws "E43" "  -3.71%  "
Set-TextValue $ws "D44" "8.64"
Set-TextValue $ws "E44" "  -1.67%  "
Set-TextValue $ws "D45" "103.23"
Set-TextValue $ws "E46" "  -2.02%  "
Set-TextValue $ws "E47" "  -2.24%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D48" "2.32"
Set-TextValue $ws "E48" "  -1.72%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D49" "1.16"
Set-TextValue $ws "E49" "  -2.85%  "
Set-TextValue $ws "D51" "2.419.24"
Set-TextValue $ws "E51" "  -1.59%  "
